$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up variant labels: strip the "T7pol_" prefix from each variant name
$ws.Range("A2").Value  = "12N"
$ws.Range("A3").Value  = "25N"
$ws.Range("A4").Value  = "WT"
$ws.Range("A5").Value  = "89R"
$ws.Range("A6").Value  = "134T"
$ws.Range("A7").Value  = "177L"
$ws.Range("A8").Value  = "225E"
$ws.Range("A9").Value  = "241W"
$ws.Range("A10").Value = "273H"

# Rename the second column header
$ws.Range("B1").Value = "fitness"

# Update the saved selection to reflect where the user left off
$ws.Range("B2").Select()
